# Generate Report for handoff
# Updates the localization-status workbook to reflect that the
# a7db943c-e265-4a52-b400-7be6b592f3d5.md file has now been handed off
# (zh-cn and de-de xlf files were generated).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# 1. "Handoff failed" -> "Not yet handed off" everywhere it is used
#    (Overview!B2, Overview!C2, zh-cn!B2, de-de!B2)
# ---------------------------------------------------------------------
$overview.Range("B2").Value = "Not yet handed off"
$overview.Range("C2").Value = "Not yet handed off"
$zhcn.Range("B2").Value = "Not yet handed off"
$dede.Range("B2").Value = "Not yet handed off"

# Helper colors/consts reused for the hyperlink-style font used elsewhere
# in the workbook (matches the existing custom "HyperLink" cell style:
# single underline, font color FF6495ED).
$hlColor = 15570276   # BGR encoding of RGB 6495ED
$hlUnderline = 2      # xlUnderlineStyleSingle

# ---------------------------------------------------------------------
# 2. zh-cn sheet, row 2 (a7db943c-...-f3d5.md)
# ---------------------------------------------------------------------
$zhcnFile = "a7db943c-e265-4a52-b400-7be6b592f3d5.f86aaf9df1a55da1992cd8a7618fb21db2b21691.zh-cn.xlf"
$zhcnUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/f86aaf9df1a55da1992cd8a7618fb21db2b21691/$zhcnFile"

$zhcn.Range("C2").Value = $zhcnFile
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), $zhcnUrl, [System.Type]::Missing, [System.Type]::Missing, $zhcnFile) | Out-Null
$zhcn.Range("C2").Font.Underline = $hlUnderline
$zhcn.Range("C2").Font.Color = $hlColor

$zhcn.Range("D2").Value = "2016-01-08 12:11:51"
$zhcn.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("H2").Value = "Include"

# ---------------------------------------------------------------------
# 3. de-de sheet, row 2 (a7db943c-...-f3d5.md)
# ---------------------------------------------------------------------
$dedeFile = "a7db943c-e265-4a52-b400-7be6b592f3d5.f86aaf9df1a55da1992cd8a7618fb21db2b21691.de-de.xlf"
$dedeUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/f86aaf9df1a55da1992cd8a7618fb21db2b21691/$dedeFile"

$dede.Range("C2").Value = $dedeFile
$dede.Hyperlinks.Add($dede.Range("C2"), $dedeUrl, [System.Type]::Missing, [System.Type]::Missing, $dedeFile) | Out-Null
$dede.Range("C2").Font.Underline = $hlUnderline
$dede.Range("C2").Font.Color = $hlColor

$dede.Range("D2").Value = "2016-01-08 12:12:05"
$dede.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("H2").Value = "Include"

Write-Host "Report generated for handoff"
